$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "70.806.29"
$ws.Cells.Item(2, 5).Value = "  +2.24%  "
$ws.Cells.Item(3, 4).Value = "3.806.24"
$ws.Cells.Item(3, 5).Value = "  +0.76%  "
$ws.Cells.Item(4, 5).Value = "  +0.08%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "699.25"
$ws.Cells.Item(5, 5).Value = "  +10.95%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "173.23"
$ws.Cells.Item(7, 4).Value = "3.806.24"
$ws.Cells.Item(7, 5).Value = "  +0.79%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "1.00"
$ws.Cells.Item(8, 5).Value = "  -0.06%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.525"
$ws.Cells.Item(9, 5).Value = "  +0.71%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.162"
$ws.Cells.Item(10, 5).Value = "  +1.72%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "7.39"
$ws.Cells.Item(11, 5).Value = "  +9.06%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.462"
$ws.Cells.Item(12, 5).Value = "  +0.50%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.0000250"
$ws.Cells.Item(13, 5).Value = "  +2.31%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "36.30"
$ws.Cells.Item(14, 5).Value = "  +3.69%  "
$ws.Cells.Item(15, 4).Value = "4.447.21"
$ws.Cells.Item(15, 5).Value = "  +0.82%  "
$ws.Cells.Item(16, 4).Value = "3.813.93"
$ws.Cells.Item(16, 5).Value = "  +1.03%  "
$ws.Cells.Item(17, 4).Value = "70.800.96"
$ws.Cells.Item(17, 5).Value = "  +2.31%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "17.77"
$ws.Cells.Item(18, 5).Value = "  +1.12%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "7.20"
$ws.Cells.Item(19, 5).Value = "  +2.57%  "
$ws.Cells.Item(20, 5).Value = "  +0.52%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "11.35"
$ws.Cells.Item(21, 5).Value = "  +18.65%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "477.67"
$ws.Cells.Item(22, 5).Value = "  +2.77%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.715"
$ws.Cells.Item(23, 5).Value = "  +1.00%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "83.77"
$ws.Cells.Item(24, 5).Value = "  +0.78%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.0000145"
$ws.Cells.Item(25, 5).Value = "  -0.51%  "
$ws.Cells.Item(26, 5).Value = "  +3.24%  "
$ws.Cells.Item(27, 2).Value = "Fetch.AI"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "2.16"
$ws.Cells.Item(27, 5).Value = "  +0.51%  "
$ws.Cells.Item(28, 2).Value = "RenderToken"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "10.39"
$ws.Cells.Item(28, 5).Value = "  +3.92%  "
$ws.Cells.Item(29, 4).Value = "3.958.20"
$ws.Cells.Item(29, 5).Value = "  +0.81%  "
$ws.Cells.Item(30, 2).Value = "Dai"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "1.00"
$ws.Cells.Item(30, 5).Value = "  -0.04%  "
$ws.Cells.Item(31, 2).Value = "PancakeSwap"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "3.13"
$ws.Cells.Item(31, 5).Value = "  +16.53%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "2.31"
$ws.Cells.Item(32, 5).Value = "  +2.11%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "7.49"
$ws.Cells.Item(33, 5).Value = "  +5.17%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "29.61"
$ws.Cells.Item(34, 5).Value = "  +3.25%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.176"
$ws.Cells.Item(35, 5).Value = "  +1.97%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "9.21"
$ws.Cells.Item(36, 5).Value = "  +2.51%  "
$ws.Cells.Item(37, 2).Value = "Binance-PegBSC-USD"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "1.00"
$ws.Cells.Item(37, 5).Value = "  +0.21%  "
$ws.Cells.Item(38, 2).Value = "RenzoRestakedETH"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Cells.Item(38, 4).Value = "3.758.85"
$ws.Cells.Item(38, 5).Value = "  +0.79%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.102"
$ws.Cells.Item(39, 5).Value = "  +0.96%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "3.50"
$ws.Cells.Item(40, 5).Value = "  +5.39%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "5.98"
$ws.Cells.Item(41, 5).Value = "  +2.83%  "
$ws.Cells.Item(42, 2).Value = "FLOKI"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.000336"
$ws.Cells.Item(42, 5).Value = "  +24.83%  "
$ws.Cells.Item(43, 2).Value = "Stacks"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "2.20"
$ws.Cells.Item(43, 5).Value = "  +13.36%  "
$ws.Cells.Item(44, 2).Value = "Mantle"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.967"
$ws.Cells.Item(44, 5).Value = "  +0.69%  "
$ws.Cells.Item(45, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(45, 5).Value = "  +0.25%  "
$ws.Cells.Item(46, 2).Value = "USDe"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "1.00"
$ws.Cells.Item(46, 5).Value = "  +0.00%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "49.52"
$ws.Cells.Item(47, 5).Value = "  +5.83%  "
$ws.Cells.Item(48, 2).Value = "Arweave"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "45.23"
$ws.Cells.Item(48, 5).Value = "  +4.38%  "
$ws.Cells.Item(49, 2).Value = "Monero"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "160.41"
$ws.Cells.Item(49, 5).Value = "  +1.53%  "
$ws.Cells.Item(50, 2).Value = "ONDO"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "1.43"
$ws.Cells.Item(50, 5).Value = "  -0.52%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.300"
$ws.Cells.Item(51, 5).Value = "  +1.82%  "
